$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OECD Data")

$ws.Range("I2").Value = 89275200000
$ws.Range("I3").Value = 91200800000
$ws.Range("I4").Value = 91268300000
$ws.Range("I5").Value = 91343600000
$ws.Range("I6").Value = 93604600000
$ws.Range("I7").Value = 97092000000
